$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 29-41: item numbers 28-40 (col A) + Danish "som var" sentences (col B)
$sentences = @(
    "Manden bar på kurven med varerne, som var",
    "Moren tændte lysene på juletræet, som var",
    "Direktøren modtog beskeden om leverancerne, som var",
    "Cyklisten vrissede over grenene på cykelstien, som var",
    "Den studerende læste bøgerne om emnet, som var",
    "Ministeren forklarede hensigten med reformerne, som var",
    "Arkæologen fremviste stenen med inskriptionerne, som var",
    "Kæresten medbragte kassen med bøgerne, som var",
    "Brandmanden slukkede branden i bygningerne, som var",
    "Pædagogen roste tegningen af kaninerne, som var",
    "Præsidenten talte om problemerne med inflationen, som var",
    "Formanden deltog i mødet om nedskæringerne, som var",
    "Sælgeren fremviste varerne fra firmaet, som var"
)

$startRow = 29
for ($i = 0; $i -lt $sentences.Length; $i++) {
    $row = $startRow + $i
    $itemNum = 28 + $i

    $ws.Cells.Item($row, 1).Value = $itemNum
    $ws.Cells.Item($row, 2).Value = $sentences[$i]

    # Alternate the shaded-row style on column A only, matching the existing
    # pattern above (even display rows -> style of A25; odd -> style of A26)
    if (($itemNum % 2) -eq 0) {
        $ws.Range("A25").Copy()
    } else {
        $ws.Range("A26").Copy()
    }
    $ws.Range("A" + $row).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the sheet's view state (scroll position + active selection)
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("B41").Select()
